$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 134.83333
$ws.Range("I53").Value = 145.625
$ws.Range("K53").Value = 145.625
$ws.Range("M53").Value = 491.375

$ws.Range("H62").Value = 13081.667
$ws.Range("I62").Value = 2650
$ws.Range("J62").Value = 18297.5
$ws.Range("K62").Value = 2650
$ws.Range("L62").Value = 18297.5
$ws.Range("M62").Value = -2026
$ws.Range("N62").Value = -19545.5

$ws.Range("H65").Value = 13081.667
$ws.Range("I65").Value = 2650
$ws.Range("J65").Value = 18297.5
$ws.Range("K65").Value = 13250
$ws.Range("L65").Value = 91487.5
$ws.Range("M65").Value = -10130
$ws.Range("N65").Value = -97727.5

$ws.Range("H127").Value = 849.3182
$ws.Range("I127").Value = 388.9
$ws.Range("J127").Value = 1233
$ws.Range("K127").Value = 1166.7
$ws.Range("L127").Value = 3699
$ws.Range("M127").Value = 3793.3
$ws.Range("N127").Value = -13619

$ws.Range("H135").Value = 15626238
$ws.Range("I135").Value = 23810602
$ws.Range("J135").Value = 1542.909
$ws.Range("K135").Value = 214295418
$ws.Range("L135").Value = 13886.181
$ws.Range("M135").Value = -214292883
$ws.Range("N135").Value = -18956.181

$ws.Range("H137").Value = 1635573.5
$ws.Range("I137").Value = 2779045.5
$ws.Range("J137").Value = 2041.9048
$ws.Range("K137").Value = 8337136.5
$ws.Range("L137").Value = 6125.7144
$ws.Range("M137").Value = -8334586.5
$ws.Range("N137").Value = -11225.7144

$ws.Range("H138").Value = 4052.9507
$ws.Range("I138").Value = 3611.9
$ws.Range("J138").Value = 4197.5576
$ws.Range("K138").Value = 10835.7
$ws.Range("L138").Value = 12592.6728
$ws.Range("M138").Value = -5695.700000000001
$ws.Range("N138").Value = -22872.6728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3495.63
$ws.Range("I32").Value = 3510.7373
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 3510.7373
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -3223.7373
$ws.Range("N32").Value = -2574

$ws.Range("H61").Value = 15155010
$ws.Range("I61").Value = 27779558
$ws.Range("J61").Value = 5551.4
$ws.Range("K61").Value = 27779558
$ws.Range("L61").Value = 5551.4
$ws.Range("M61").Value = -27779346
$ws.Range("N61").Value = -5975.4

$ws.Range("H63").Value = 104410.89
$ws.Range("J63").Value = 4962.25
$ws.Range("L63").Value = 4962.25
$ws.Range("N63").Value = -6334.25

$ws.Range("H66").Value = 104410.89
$ws.Range("J66").Value = 4962.25
$ws.Range("L66").Value = 24811.25
$ws.Range("N66").Value = -31675.25

$ws.Range("H74").Value = 13515342
$ws.Range("I74").Value = 917.3333
$ws.Range("J74").Value = 38465052
$ws.Range("K74").Value = 917.3333
$ws.Range("L74").Value = 38465052
$ws.Range("M74").Value = -43.33330000000001
$ws.Range("N74").Value = -38466800

$ws.Range("H77").Value = 13515342
$ws.Range("I77").Value = 917.3333
$ws.Range("J77").Value = 38465052
$ws.Range("K77").Value = 4586.6665
$ws.Range("L77").Value = 192325260
$ws.Range("M77").Value = -218.6665000000003
$ws.Range("N77").Value = -192333996

$ws.Range("H132").Value = 1101533.5
$ws.Range("I132").Value = 1823
$ws.Range("J132").Value = 5133805.5
$ws.Range("K132").Value = 5469
$ws.Range("L132").Value = 15401416.5
$ws.Range("M132").Value = -2939
$ws.Range("N132").Value = -15406476.5

$ws.Range("H136").Value = 15155010
$ws.Range("I136").Value = 27779558
$ws.Range("J136").Value = 5551.4
$ws.Range("K136").Value = 83338674
$ws.Range("L136").Value = 16654.2
$ws.Range("M136").Value = -83336124
$ws.Range("N136").Value = -21754.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 95657
$ws.Range("J103").Value = 95657
$ws.Range("L103").Value = 95657
$ws.Range("N103").Value = -98001

$ws.Range("H107").Value = 85525.586
$ws.Range("I107").Value = 102030.7
$ws.Range("K107").Value = 102030.7
$ws.Range("M107").Value = -100110.7

$ws.Range("H134").Value = 2606.1
$ws.Range("I134").Value = 2517.2
$ws.Range("J134").Value = 2872.8
$ws.Range("K134").Value = 7551.599999999999
$ws.Range("L134").Value = 8618.400000000001
$ws.Range("M134").Value = -5016.599999999999
$ws.Range("N134").Value = -13688.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5212457
$ws.Range("I134").Value = 6948692
$ws.Range("K134").Value = 20846076
$ws.Range("M134").Value = -20843541

$ws.Range("H140").Value = 64262.5
$ws.Range("J140").Value = 64262.5
$ws.Range("L140").Value = 64262.5
$ws.Range("N140").Value = -74622.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 550
$ws.Range("I98").Value = 400
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 1200
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 298
$ws.Range("N98").Value = -5996

$ws.Range("H122").Value = 2956.9216
$ws.Range("I122").Value = 505.6842
$ws.Range("J122").Value = 4412.3438
$ws.Range("K122").Value = 4551.1578
$ws.Range("L122").Value = 39711.0942
$ws.Range("M122").Value = -2101.1578
$ws.Range("N122").Value = -44611.0942

$ws.Range("H132").Value = 2685.5
$ws.Range("I132").Value = 2168
$ws.Range("J132").Value = 3203
$ws.Range("K132").Value = 19512
$ws.Range("L132").Value = 28827
$ws.Range("M132").Value = -16982
$ws.Range("N132").Value = -33887

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 45461996
$ws.Range("I132").Value = 62508364
$ws.Range("K132").Value = 187525092
$ws.Range("M132").Value = -187522562

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5200.5557
$ws.Range("I61").Value = 2825
$ws.Range("K61").Value = 2825
$ws.Range("M61").Value = -2623

$ws.Range("H113").Value = 5200.5557
$ws.Range("I113").Value = 2825
$ws.Range("K113").Value = 2825
$ws.Range("M113").Value = -655

$ws.Range("H122").Value = 4906.5776
$ws.Range("I122").Value = 3514.0715
$ws.Range("J122").Value = 5535.4517
$ws.Range("K122").Value = 10542.2145
$ws.Range("L122").Value = 16606.3551
$ws.Range("M122").Value = -8092.2145
$ws.Range("N122").Value = -21506.3551

$ws.Range("H136").Value = 15154634
$ws.Range("I136").Value = 2326
$ws.Range("J136").Value = 23813096
$ws.Range("K136").Value = 6978
$ws.Range("L136").Value = 71439288
$ws.Range("M136").Value = -4428
$ws.Range("N136").Value = -71444388

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4783858
$ws.Range("I132").Value = 2443.6667
$ws.Range("J132").Value = 15353301
$ws.Range("K132").Value = 7331.000100000001
$ws.Range("L132").Value = 46059903
$ws.Range("M132").Value = -4801.000100000001
$ws.Range("N132").Value = -46064963
